$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple D/E value updates (price + volume % change) ---
$ws.Cells.Item(2, 4).Value = '19.865.01'
$ws.Cells.Item(2, 5).Value = '  -8.41%  '

$ws.Cells.Item(3, 4).Value = '1.402.88'
$ws.Cells.Item(3, 5).Value = '  -8.80%  '

$ws.Cells.Item(4, 5).Value = '  +0.35%  '

$c3 = $ws.Cells.Item(5, 4)
$c3.Value = "'1.003"
$c3.ClearFormats()
$ws.Cells.Item(5, 5).Value = '  +0.39%  '

$ws.Cells.Item(6, 5).Value = '  -5.71%  '

$c4 = $ws.Cells.Item(7, 4)
$c4.Value = "'0.3679"
$c4.ClearFormats()
$ws.Cells.Item(7, 5).Value = '  -6.50%  '

$c5 = $ws.Cells.Item(8, 4)
$c5.Value = "'0.3108"
$c5.ClearFormats()
$ws.Cells.Item(8, 5).Value = '  -2.78%  '

$c6 = $ws.Cells.Item(9, 4)
$c6.Value = "'39.47"
$c6.ClearFormats()
$ws.Cells.Item(9, 5).Value = '  -7.98%  '

$c7 = $ws.Cells.Item(10, 4)
$c7.Value = "'1.004"
$c7.ClearFormats()
$ws.Cells.Item(10, 5).Value = '  -6.93%  '

$c8 = $ws.Cells.Item(11, 4)
$c8.Value = "'0.06471"
$c8.ClearFormats()
$ws.Cells.Item(11, 5).Value = '  -10.09%  '

$ws.Cells.Item(12, 5).Value = '  +0.33%  '

$c9 = $ws.Cells.Item(13, 4)
$c9.Value = "'5.437"
$c9.ClearFormats()
$ws.Cells.Item(13, 5).Value = '  -5.56%  '

$c10 = $ws.Cells.Item(14, 4)
$c10.Value = "'17.35"
$c10.ClearFormats()
$ws.Cells.Item(14, 5).Value = '  -5.95%  '

$c11 = $ws.Cells.Item(15, 4)
$c11.Value = "'6.121"
$c11.ClearFormats()
$ws.Cells.Item(15, 5).Value = '  -7.86%  '

$ws.Cells.Item(16, 4).Value = '1.400.10'
$ws.Cells.Item(16, 5).Value = '  -9.47%  '

$c13 = $ws.Cells.Item(17, 4)
$c13.Value = "'0.00001007"
$c13.ClearFormats()
$ws.Cells.Item(17, 5).Value = '  -8.21%  '

$c14 = $ws.Cells.Item(18, 4)
$c14.Value = "'0.05683"
$c14.ClearFormats()
$ws.Cells.Item(18, 5).Value = '  -14.01%  '

$ws.Cells.Item(19, 5).Value = '  +0.50%  '

$c15 = $ws.Cells.Item(20, 4)
$c15.Value = "'70.08"
$c15.ClearFormats()
$ws.Cells.Item(20, 5).Value = '  -16.89%  '

$c16 = $ws.Cells.Item(21, 4)
$c16.Value = "'5.558"
$c16.ClearFormats()
$ws.Cells.Item(21, 5).Value = '  -9.83%  '

$c17 = $ws.Cells.Item(22, 4)
$c17.Value = "'14.65"
$c17.ClearFormats()
$ws.Cells.Item(22, 5).Value = '  -6.11%  '

$c18 = $ws.Cells.Item(23, 4)
$c18.Value = "'11.02"
$c18.ClearFormats()
$ws.Cells.Item(23, 5).Value = '  +1.31%  '

$c19 = $ws.Cells.Item(24, 4)
$c19.Value = "'2.280"
$c19.ClearFormats()
$ws.Cells.Item(24, 5).Value = '  -3.55%  '

$ws.Cells.Item(25, 4).Value = '19.888.09'
$ws.Cells.Item(25, 5).Value = '  -8.34%  '

$c21 = $ws.Cells.Item(26, 4)
$c21.Value = "'2.216"
$c21.ClearFormats()
$ws.Cells.Item(26, 5).Value = '  -7.74%  '

$c22 = $ws.Cells.Item(27, 4)
$c22.Value = "'135.04"
$c22.ClearFormats()
$ws.Cells.Item(27, 5).Value = '  -11.16%  '

$c23 = $ws.Cells.Item(28, 4)
$c23.Value = "'16.86"
$c23.ClearFormats()
$ws.Cells.Item(28, 5).Value = '  -9.00%  '

$ws.Cells.Item(29, 4).Value = '1.557.55'
$ws.Cells.Item(29, 5).Value = '  -9.48%  '

$c25 = $ws.Cells.Item(30, 4)
$c25.Value = "'109.15"
$c25.ClearFormats()
$ws.Cells.Item(30, 5).Value = '  -7.41%  '

$c26 = $ws.Cells.Item(31, 4)
$c26.Value = "'4.090"
$c26.ClearFormats()
$ws.Cells.Item(31, 5).Value = '  -15.92%  '

$c27 = $ws.Cells.Item(32, 4)
$c27.Value = "'5.258"
$c27.ClearFormats()
$ws.Cells.Item(32, 5).Value = '  -14.26%  '

$c28 = $ws.Cells.Item(33, 4)
$c28.Value = "'0.8077"
$c28.ClearFormats()
$ws.Cells.Item(33, 5).Value = '  -16.86%  '

$c29 = $ws.Cells.Item(34, 4)
$c29.Value = "'0.07640"
$c29.ClearFormats()
$ws.Cells.Item(34, 5).Value = '  -6.19%  '

$c30 = $ws.Cells.Item(35, 4)
$c30.Value = "'8.301"
$c30.ClearFormats()
$ws.Cells.Item(35, 5).Value = '  -3.56%  '

$c31 = $ws.Cells.Item(38, 4)
$c31.Value = "'4.802"
$c31.ClearFormats()
$ws.Cells.Item(38, 5).Value = '  -7.83%  '

$ws.Cells.Item(39, 5).Value = '  +0.38%  '

$c32 = $ws.Cells.Item(40, 4)
$c32.Value = "'0.02058"
$c32.ClearFormats()
$ws.Cells.Item(40, 5).Value = '  -8.19%  '

$c33 = $ws.Cells.Item(41, 4)
$c33.Value = "'0.1892"
$c33.ClearFormats()
$ws.Cells.Item(41, 5).Value = '  -7.76%  '

$ws.Cells.Item(42, 5).Value = '  -9.25%  '

$ws.Cells.Item(43, 5).Value = '  -8.03%  '

$c34 = $ws.Cells.Item(46, 4)
$c34.Value = "'12.15"
$c34.ClearFormats()
$ws.Cells.Item(46, 5).Value = '  -8.76%  '

$c35 = $ws.Cells.Item(47, 4)
$c35.Value = "'0.5072"
$c35.ClearFormats()
$ws.Cells.Item(47, 5).Value = '  -9.25%  '

$c36 = $ws.Cells.Item(48, 4)
$c36.Value = "'111.83"
$c36.ClearFormats()
$ws.Cells.Item(48, 5).Value = '  -4.32%  '

$c37 = $ws.Cells.Item(49, 4)
$c37.Value = "'1.754"
$c37.ClearFormats()
$ws.Cells.Item(49, 5).Value = '  -7.65%  '

# --- Row swaps: coin ranking order changed, B/C/D/E updated together ---
$ws.Cells.Item(36, 2).Value = 'Hedera'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c38 = $ws.Cells.Item(36, 4)
$c38.Value = "'0.05768"
$c38.ClearFormats()
$ws.Cells.Item(36, 5).Value = '  -4.28%  '

$ws.Cells.Item(37, 2).Value = 'WEMIXTOKEN'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$c39 = $ws.Cells.Item(37, 4)
$c39.Value = "'1.427"
$c39.ClearFormats()
$ws.Cells.Item(37, 5).Value = '  -4.37%  '

$ws.Cells.Item(44, 2).Value = 'TheSandbox'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$c40 = $ws.Cells.Item(44, 4)
$c40.Value = "'0.5240"
$c40.ClearFormats()
$ws.Cells.Item(44, 5).Value = '  -10.30%  '

$ws.Cells.Item(45, 2).Value = 'PancakeSwap'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$c41 = $ws.Cells.Item(45, 4)
$c41.Value = "'3.501"
$c41.ClearFormats()
$ws.Cells.Item(45, 5).Value = '  -6.12%  '

$ws.Cells.Item(50, 2).Value = 'EOS'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$c42 = $ws.Cells.Item(50, 4)
$c42.Value = "'1.035"
$c42.ClearFormats()
$ws.Cells.Item(50, 5).Value = '  -11.49%  '

$ws.Cells.Item(51, 2).Value = 'PaxDollar'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$c43 = $ws.Cells.Item(51, 4)
$c43.Value = "'1.004"
$c43.ClearFormats()
$ws.Cells.Item(51, 5).Value = '  +0.32%  '

